# Weekly NYPD CompStat refresh (40th Precinct) - new crime data collected.
#
# Updates the report header (volume/report-week labels) and the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year crime-complaint figures
# (and their derived % Chg columns) for rows 14-29 of the CompStat sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header -------------------------------------------------------
$ws.Range("A8").Value2 = "Volume 30   Number  40"
$ws.Range("C9").Value2 = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Crime-complaint table -------------------------------------------------
# Columns C..N map to:
#   C = WTD 2023, D = WTD 2022, E = WTD %Chg,
#   F = 28-Day 2023, G = 28-Day 2022, H = 28-Day %Chg,
#   I = YTD 2023, J = YTD 2022, K = YTD %Chg,
#   L = 2-Year %Chg, M = 13-Year %Chg, N = 30-Year %Chg
# A value of $null means that cell is unchanged this week.
$rowData = [ordered]@{
    14 = @($null, $null, $null, 1, 3, -66.666666666666, $null, 12, -33.333333333333, -42.857142857142, $null, -86.206896551724)
    15 = @($null, $null, $null, $null, 4, -25, 35, 26, 34.615384615384, 59.090909090909, 105.882352941176, -37.5)
    16 = @(20, 8, 150, 74, 35, 111.428571428571, 521, 441, 18.140589569161, 65.923566878980, 46.348314606741, -64.339493497604)
    17 = @(26, 12, 116.666666666667, 96, 74, 29.729729729729, 788, 658, 19.756838905775, 45.925925925925, 132.448377581121, -9.321058688147)
    18 = @(6, 7, -14.285714285714, $null, 38, -15.789473684210, 247, 278, -11.151079136690, 85.714285714285, 57.324840764331, -76.829268292682)
    19 = @(11, 18, -38.888888888888, 50, $null, -26.470588235294, 550, 586, -6.143344709897, 0.917431192660, 66.163141993957, -8.789386401326)
    20 = @(7, 5, 40, $null, 21, 19.047619047619, 280, 244, 14.754098360655, 101.438848920863, 191.666666666667, -46.969696969697)
    21 = @($null, 52, 36.538461538461, 281, 243, 15.637860082304, 2429, 2245, 8.195991091314, 42.296426479203, 85.561497326203, -47.662141779788)
    22 = @(2, 1, 100, 5, $null, 0, 50, 72, -30.555555555555, -7.407407407407, 4.166666666666, $null)
    23 = @(6, 5, 20, 27, 27, 0, 378, 309, 22.330097087378, 72.602739726027, 60.169491525423, $null)
    24 = @(31, 36, -13.888888888888, 142, 233, -39.055793991416, 1268, 1425, -11.017543859649, 23.346303501945, 19.735599622285, $null)
    25 = @(25, 8, 212.5, 94, 91, 3.296703296703, 873, 797, 9.535759096612, 23.654390934844, -1.689189189189, $null)
    26 = @($null, $null, $null, 4, $null, -33.333333333333, 49, 43, 13.953488372093, -10.909090909090, $null, $null)
    27 = @(4, 1, 300, 13, 5, 160, 92, 68, 35.294117647058, 41.538461538461, $null, $null)
    28 = @($null, 1, $null, 4, 4, $null, $null, 51, -49.019607843137, -52.727272727272, -48, -85.310734463276)
    29 = @($null, 1, $null, 3, 4, -25, $null, 42, -42.857142857142, -50, -42.857142857142, -84.810126582278)
}

$columns = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $newValue = $values[$i]
        if ($null -ne $newValue) {
            $ws.Range("$($columns[$i])$row").Value2 = $newValue
        }
    }
}

# A couple of cells that previously had no incidents ("***.*" / N/A text)
# now carry real figures, so they need to switch from the text format back
# to the normal numeric / %-change number formats used elsewhere in the table.
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
